$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unneeded footnote/contact block (rows 21-31) contents.
$ws.Range("A21:G31").ClearContents()

# Restore the selection the author left behind in the saved file.
$ws.Range("A19:G33").Select()
